# doc/journal_de_travail.xlsx — log three new "Réalisation du modèle" work
# sessions (model training / export notebook work) into the Journal sheet,
# rows 24-26, which were previously blank placeholder rows (only the
# carried-down D-column duration formula existed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$ws.Activate()

# Row 24: 2024-07-27, 14:00 -> 18:00, "Réalisation du modèle"
$ws.Range("A24").Value = 45500
$ws.Range("B24").Value = 0.58333333333333337
$ws.Range("C24").Value = 0.75
$ws.Range("E24").Value = "Réalisation du modèle"

# Row 25: 2024-08-08, 13:47 -> 17:25, "Réalisation du modèle"
$ws.Range("A25").Value = 45512
$ws.Range("B25").Value = 0.57430555555555551
$ws.Range("C25").Value = 0.72569444444444453
$ws.Range("E25").Value = "Réalisation du modèle"

# Row 26: 2024-08-09, 13:38 -> 18:00, "Réalisation du modèle"
$ws.Range("A26").Value = 45513
$ws.Range("B26").Value = 0.56805555555555554
$ws.Range("C26").Value = 0.75
$ws.Range("E26").Value = "Réalisation du modèle"

# Recalculate so the dependent D-column durations and the H/I SUMIF summary
# table (including the "Total" in I7) pick up the three new rows.
$excel.CalculateFull()

# Leave the sheet scrolled/selected where the user ended up after typing
# the new rows.
[void]$ws.Range("E27").Select()
